$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(1, 2).Value = "bank"
$ws.Cells.Item(1, 3).Value = "deposit_type"
$ws.Cells.Item(1, 4).Value = "currency"
$ws.Cells.Item(1, 5).Value = "owner"
$ws.Cells.Item(1, 6).Value = "total"
$ws.Cells.Item(1, 7).Value = "property_category"
$ws.Cells.Item(1, 8).Value = "category"
$ws.Cells.Item(1, 9).Value = "date"
$ws.Cells.Item(1, 10).Value = "legislator_name"
$ws.Cells.Item(1, 11).Value = "legislator_id"
$ws.Cells.Item(1, 12).Value = "source_file"
$ws.Cells.Item(1, 13).Value = "index"
$ws.Cells.Item(2, 1).Value = 44
$ws.Cells.Item(2, 2).Value = "臺灣銀行城中分行"
$ws.Cells.Item(2, 3).Value = "活期存款"
$ws.Cells.Item(2, 4).Value = "美金"
$ws.Cells.Item(2, 5).Value = "邱議瑩"
$ws.Cells.Item(2, 6).Value = 58320
$ws.Cells.Item(2, 7).Value = "deposit"
$ws.Cells.Item(2, 8).Value = "normal"
$ws.Cells.Item(2, 9).Value = "2012-04-18"
$ws.Cells.Item(2, 10).Value = "邱議瑩"
$ws.Cells.Item(2, 11).Value = 913
$ws.Cells.Item(2, 12).Value = "tmped121"
$ws.Cells.Item(2, 13).Value = 44
$ws.Cells.Item(3, 1).Value = 45
$ws.Cells.Item(3, 2).Value = "臺灣銀行群賢分行"
$ws.Cells.Item(3, 3).Value = "活期存款"
$ws.Cells.Item(3, 4).Value = "新臺幣"
$ws.Cells.Item(3, 5).Value = "邱議瑩"
$ws.Cells.Item(3, 6).Value = 7200944
$ws.Cells.Item(3, 7).Value = "deposit"
$ws.Cells.Item(3, 8).Value = "normal"
$ws.Cells.Item(3, 9).Value = "2012-04-18"
$ws.Cells.Item(3, 10).Value = "邱議瑩"
$ws.Cells.Item(3, 11).Value = 913
$ws.Cells.Item(3, 12).Value = "tmped121"
$ws.Cells.Item(3, 13).Value = 45
$ws.Cells.Item(4, 1).Value = 46
$ws.Cells.Item(4, 2).Value = "臺灣新光商業銀行屏東分行"
$ws.Cells.Item(4, 3).Value = "定期儲蓄存款"
$ws.Cells.Item(4, 4).Value = "新臺幣"
$ws.Cells.Item(4, 5).Value = "邱議瑩"
$ws.Cells.Item(4, 6).Value = 300000
$ws.Cells.Item(4, 7).Value = "deposit"
$ws.Cells.Item(4, 8).Value = "normal"
$ws.Cells.Item(4, 9).Value = "2012-04-18"
$ws.Cells.Item(4, 10).Value = "邱議瑩"
$ws.Cells.Item(4, 11).Value = 913
$ws.Cells.Item(4, 12).Value = "tmped121"
$ws.Cells.Item(4, 13).Value = 46
$ws.Cells.Item(5, 1).Value = 47
$ws.Cells.Item(5, 2).Value = "臺灣新光商業銀行屏東分行"
$ws.Cells.Item(5, 3).Value = "活期存款"
$ws.Cells.Item(5, 4).Value = "新臺幣"
$ws.Cells.Item(5, 5).Value = "邱議瑩"
$ws.Cells.Item(5, 6).Value = 1800000
$ws.Cells.Item(5, 7).Value = "deposit"
$ws.Cells.Item(5, 8).Value = "normal"
$ws.Cells.Item(5, 9).Value = "2012-04-18"
$ws.Cells.Item(5, 10).Value = "邱議瑩"
$ws.Cells.Item(5, 11).Value = 913
$ws.Cells.Item(5, 12).Value = "tmped121"
$ws.Cells.Item(5, 13).Value = 47
$ws.Cells.Item(6, 1).Value = 48
$ws.Cells.Item(6, 2).Value = "彰化商業銀行屏東分行"
$ws.Cells.Item(6, 3).Value = "活期存款"
$ws.Cells.Item(6, 4).Value = "新臺幣"
$ws.Cells.Item(6, 5).Value = "邱議瑩"
$ws.Cells.Item(6, 6).Value = 2835
$ws.Cells.Item(6, 7).Value = "deposit"
$ws.Cells.Item(6, 8).Value = "normal"
$ws.Cells.Item(6, 9).Value = "2012-04-18"
$ws.Cells.Item(6, 10).Value = "邱議瑩"
$ws.Cells.Item(6, 11).Value = 913
$ws.Cells.Item(6, 12).Value = "tmped121"
$ws.Cells.Item(6, 13).Value = 48
$ws.Cells.Item(7, 1).Value = 49
$ws.Cells.Item(7, 2).Value = "合作金庫商業銀行營業部"
$ws.Cells.Item(7, 3).Value = "活期存款"
$ws.Cells.Item(7, 4).Value = "新臺幣"
$ws.Cells.Item(7, 5).Value = "邱議瑩"
$ws.Cells.Item(7, 6).Value = 2057
$ws.Cells.Item(7, 7).Value = "deposit"
$ws.Cells.Item(7, 8).Value = "normal"
$ws.Cells.Item(7, 9).Value = "2012-04-18"
$ws.Cells.Item(7, 10).Value = "邱議瑩"
$ws.Cells.Item(7, 11).Value = 913
$ws.Cells.Item(7, 12).Value = "tmped121"
$ws.Cells.Item(7, 13).Value = 49
$ws.Cells.Item(8, 1).Value = 50
$ws.Cells.Item(8, 2).Value = "元大商業銀行營業部"
$ws.Cells.Item(8, 3).Value = "活期存款"
$ws.Cells.Item(8, 4).Value = "新臺幣"
$ws.Cells.Item(8, 5).Value = "邱議瑩"
$ws.Cells.Item(8, 6).Value = 102866
$ws.Cells.Item(8, 7).Value = "deposit"
$ws.Cells.Item(8, 8).Value = "normal"
$ws.Cells.Item(8, 9).Value = "2012-04-18"
$ws.Cells.Item(8, 10).Value = "邱議瑩"
$ws.Cells.Item(8, 11).Value = 913
$ws.Cells.Item(8, 12).Value = "tmped121"
$ws.Cells.Item(8, 13).Value = 50
$ws.Cells.Item(9, 1).Value = 51
$ws.Cells.Item(9, 2).Value = "京城商業銀行營業部"
$ws.Cells.Item(9, 3).Value = "活期存款"
$ws.Cells.Item(9, 4).Value = "新臺幣"
$ws.Cells.Item(9, 5).Value = "邱議瑩"
$ws.Cells.Item(9, 6).Value = 3364
$ws.Cells.Item(9, 7).Value = "deposit"
$ws.Cells.Item(9, 8).Value = "normal"
$ws.Cells.Item(9, 9).Value = "2012-04-18"
$ws.Cells.Item(9, 10).Value = "邱議瑩"
$ws.Cells.Item(9, 11).Value = 913
$ws.Cells.Item(9, 12).Value = "tmped121"
$ws.Cells.Item(9, 13).Value = 51
$ws.Cells.Item(10, 1).Value = 52
$ws.Cells.Item(10, 2).Value = "NationalAustraliaBankSydency"
$ws.Cells.Item(10, 3).Value = "活期存款"
$ws.Cells.Item(10, 4).Value = "澳幣"
$ws.Cells.Item(10, 5).Value = "李永得"
$ws.Cells.Item(10, 6).Value = 216840
$ws.Cells.Item(10, 7).Value = "deposit"
$ws.Cells.Item(10, 8).Value = "normal"
$ws.Cells.Item(10, 9).Value = "2012-04-18"
$ws.Cells.Item(10, 10).Value = "邱議瑩"
$ws.Cells.Item(10, 11).Value = 913
$ws.Cells.Item(10, 12).Value = "tmped121"
$ws.Cells.Item(10, 13).Value = 52
$ws.Cells.Item(11, 1).Value = 53
$ws.Cells.Item(11, 2).Value = "高雄銀行市府分行"
$ws.Cells.Item(11, 3).Value = "活期存款"
$ws.Cells.Item(11, 4).Value = "新臺幣"
$ws.Cells.Item(11, 5).Value = "李永得"
$ws.Cells.Item(11, 6).Value = 5598972
$ws.Cells.Item(11, 7).Value = "deposit"
$ws.Cells.Item(11, 8).Value = "normal"
$ws.Cells.Item(11, 9).Value = "2012-04-18"
$ws.Cells.Item(11, 10).Value = "邱議瑩"
$ws.Cells.Item(11, 11).Value = 913
$ws.Cells.Item(11, 12).Value = "tmped121"
$ws.Cells.Item(11, 13).Value = 53
$ws.Cells.Item(12, 1).Value = 54
$ws.Cells.Item(12, 2).Value = "高雄銀行市府分行"
$ws.Cells.Item(12, 3).Value = "公教優惠儲蓄存款"
$ws.Cells.Item(12, 4).Value = "新臺幣"
$ws.Cells.Item(12, 5).Value = "李永得"
$ws.Cells.Item(12, 6).Value = 437987
$ws.Cells.Item(12, 7).Value = "deposit"
$ws.Cells.Item(12, 8).Value = "normal"
$ws.Cells.Item(12, 9).Value = "2012-04-18"
$ws.Cells.Item(12, 10).Value = "邱議瑩"
$ws.Cells.Item(12, 11).Value = 913
$ws.Cells.Item(12, 12).Value = "tmped121"
$ws.Cells.Item(12, 13).Value = 54
$ws.Cells.Item(13, 1).Value = 55
$ws.Cells.Item(13, 2).Value = "高雄銀行市府分行"
$ws.Cells.Item(13, 3).Value = "活期存款"
$ws.Cells.Item(13, 4).Value = "新加坡幣"
$ws.Cells.Item(13, 5).Value = "李永得"
$ws.Cells.Item(13, 6).Value = 853110
$ws.Cells.Item(13, 7).Value = "deposit"
$ws.Cells.Item(13, 8).Value = "normal"
$ws.Cells.Item(13, 9).Value = "2012-04-18"
$ws.Cells.Item(13, 10).Value = "邱議瑩"
$ws.Cells.Item(13, 11).Value = 913
$ws.Cells.Item(13, 12).Value = "tmped121"
$ws.Cells.Item(13, 13).Value = 55
$ws.Cells.Item(14, 1).Value = 56
$ws.Cells.Item(14, 2).Value = "臺灣銀行松山分行"
$ws.Cells.Item(14, 3).Value = "活期存款"
$ws.Cells.Item(14, 4).Value = "新臺幣"
$ws.Cells.Item(14, 5).Value = "李永得"
$ws.Cells.Item(14, 6).Value = 21028
$ws.Cells.Item(14, 7).Value = "deposit"
$ws.Cells.Item(14, 8).Value = "normal"
$ws.Cells.Item(14, 9).Value = "2012-04-18"
$ws.Cells.Item(14, 10).Value = "邱議瑩"
$ws.Cells.Item(14, 11).Value = 913
$ws.Cells.Item(14, 12).Value = "tmped121"
$ws.Cells.Item(14, 13).Value = 56
$ws.Cells.Item(15, 1).Value = 57
$ws.Cells.Item(15, 2).Value = "國泰世華商業銀行南京東路分行"
$ws.Cells.Item(15, 3).Value = "活期存款"
$ws.Cells.Item(15, 4).Value = "新臺幣"
$ws.Cells.Item(15, 5).Value = "李永得"
$ws.Cells.Item(15, 6).Value = 99657
$ws.Cells.Item(15, 7).Value = "deposit"
$ws.Cells.Item(15, 8).Value = "normal"
$ws.Cells.Item(15, 9).Value = "2012-04-18"
$ws.Cells.Item(15, 10).Value = "邱議瑩"
$ws.Cells.Item(15, 11).Value = 913
$ws.Cells.Item(15, 12).Value = "tmped121"
$ws.Cells.Item(15, 13).Value = 57
$ws.Cells.Item(16, 1).Value = 58
$ws.Cells.Item(16, 2).Value = "第一商業銀行總行"
$ws.Cells.Item(16, 3).Value = "活期存款"
$ws.Cells.Item(16, 4).Value = "新臺幣"
$ws.Cells.Item(16, 5).Value = "李永得"
$ws.Cells.Item(16, 6).Value = 1632
$ws.Cells.Item(16, 7).Value = "deposit"
$ws.Cells.Item(16, 8).Value = "normal"
$ws.Cells.Item(16, 9).Value = "2012-04-18"
$ws.Cells.Item(16, 10).Value = "邱議瑩"
$ws.Cells.Item(16, 11).Value = 913
$ws.Cells.Item(16, 12).Value = "tmped121"
$ws.Cells.Item(16, 13).Value = 58
$ws.Cells.Item(17, 1).Value = 59
$ws.Cells.Item(17, 2).Value = "臺灣土地銀行總行"
$ws.Cells.Item(17, 3).Value = "活期存款"
$ws.Cells.Item(17, 4).Value = "新臺幣"
$ws.Cells.Item(17, 5).Value = "李永得"
$ws.Cells.Item(17, 6).Value = 426
$ws.Cells.Item(17, 7).Value = "deposit"
$ws.Cells.Item(17, 8).Value = "normal"
$ws.Cells.Item(17, 9).Value = "2012-04-18"
$ws.Cells.Item(17, 10).Value = "邱議瑩"
$ws.Cells.Item(17, 11).Value = 913
$ws.Cells.Item(17, 12).Value = "tmped121"
$ws.Cells.Item(17, 13).Value = 59
$ws.Cells.Item(18, 1).Value = 60
$ws.Cells.Item(18, 2).Value = "兆豐國際商業銀行總管理處"
$ws.Cells.Item(18, 3).Value = "活期存款"
$ws.Cells.Item(18, 4).Value = "新臺幣"
$ws.Cells.Item(18, 5).Value = "李永得"
$ws.Cells.Item(18, 6).Value = 19411
$ws.Cells.Item(18, 7).Value = "deposit"
$ws.Cells.Item(18, 8).Value = "normal"
$ws.Cells.Item(18, 9).Value = "2012-04-18"
$ws.Cells.Item(18, 10).Value = "邱議瑩"
$ws.Cells.Item(18, 11).Value = 913
$ws.Cells.Item(18, 12).Value = "tmped121"
$ws.Cells.Item(18, 13).Value = 60
$ws.Cells.Item(19, 1).Value = 61
$ws.Cells.Item(19, 2).Value = "中國信託商業銀行敦北分行"
$ws.Cells.Item(19, 3).Value = "活期存款"
$ws.Cells.Item(19, 4).Value = "新臺幣"
$ws.Cells.Item(19, 5).Value = "李永得"
$ws.Cells.Item(19, 6).Value = 5530
$ws.Cells.Item(19, 7).Value = "deposit"
$ws.Cells.Item(19, 8).Value = "normal"
$ws.Cells.Item(19, 9).Value = "2012-04-18"
$ws.Cells.Item(19, 10).Value = "邱議瑩"
$ws.Cells.Item(19, 11).Value = 913
$ws.Cells.Item(19, 12).Value = "tmped121"
$ws.Cells.Item(19, 13).Value = 61
$ws.Cells.Item(20, 1).Value = 62
$ws.Cells.Item(20, 2).Value = "花旗（台灣）銀行營業部"
$ws.Cells.Item(20, 3).Value = "活期存款"
$ws.Cells.Item(20, 4).Value = "新臺幣"
$ws.Cells.Item(20, 5).Value = "李永得"
$ws.Cells.Item(20, 6).Value = 2705
$ws.Cells.Item(20, 7).Value = "deposit"
$ws.Cells.Item(20, 8).Value = "normal"
$ws.Cells.Item(20, 9).Value = "2012-04-18"
$ws.Cells.Item(20, 10).Value = "邱議瑩"
$ws.Cells.Item(20, 11).Value = 913
$ws.Cells.Item(20, 12).Value = "tmped121"
$ws.Cells.Item(20, 13).Value = 62
